$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 headers
$ws.Range("A1").Value = "EsriCareersMainURL"
$ws.Range("B1").Value = "careerMHdr"
$ws.Range("C1").Value = "JDHdr"
$ws.Range("D1").Value = "Inputdata"
$ws.Range("E1").Value = "jobCategory"
$ws.Range("F1").Value = "Location"
$ws.Range("G1").Value = "LHdr"
$ws.Range("H1").Value = "REHeader"
$ws.Range("I1").Value = "StuHdr"
$ws.Range("J1").Value = "JCitemdata"
$ws.Range("K1").Value = "LifeHdr"
$ws.Range("L1").Value = "CHeader"
$ws.Range("M1").Value = "BText"
$ws.Range("N1").Value = "FName"
$ws.Range("O1").Value = "LName"
$ws.Range("P1").Value = "EmailID"
$ws.Range("Q1").Value = "PhoneNo"
$ws.Range("R1").Value = "VConfNotific"

# Row 2 data
$ws.Range("C2").Value = "JOB OPENINGS"
$ws.Range("D2").Value = "Testing"
$ws.Range("E2").Value = "Information Technology"
$ws.Range("G2").Value = "US-CA-Redlands"
$ws.Range("H2").Value = "RECRUITING EVENTS"
$ws.Range("I2").Value = "STUDENTS AT ESRI"
$ws.Range("J2").Value = "Administrative Support"
$ws.Range("K2").Value = "LIFE @ ESRI"
$ws.Range("L2").Value = "CONTACT US"
$ws.Range("M2").Value = "Testing"
$ws.Range("N2").Value = "Test1"
$ws.Range("O2").Value = "Ltest"
$ws.Range("P2").Value = "balaji.harinath@htcindia.com"
$ws.Range("Q2").Value = 44222222222
$ws.Range("R2").Value = "Thank You!"
$ws.Range("B2").Value = "We Are Esri"

# Hyperlink on the email cell, keep the existing Hyperlink cell style index
$ws.Hyperlinks.Add($ws.Range("P2"), "mailto:balaji.harinath@htcindia.com") | Out-Null
$ws.Range("P2").Style = "Hyperlink"

# Column width adjustments
$ws.Columns(2).ColumnWidth = 11.1640625
$ws.Columns(17).ColumnWidth = 11.1640625

# Selection moves to B3
$ws.Range("B3").Select() | Out-Null
